# Adds a second worksheet ("data2") that holds:
#   - a small "verification number" table (A1:A2) used (in the original
#     authored workbook) as the source of a Power Query parameter query
#   - a "loaded query result" table (A4:B5) that mirrors what the
#     Power Query "query2" query would load back into the sheet
# together with the supporting table objects, the hidden
# "ExternalData_1" defined name Excel creates for a query-table backed
# table, and the view/selection state the diff captures.
#
# NOTE: the source workbook also gained genuine Power Query plumbing
# (xl/connections.xml, xl/queryTables/queryTable1.xml, a
# tableType="queryTable" table, and a customXml DataMashup part holding
# the M code). Those artifacts are produced by Excel's Get&Transform /
# Power Query engine, which has no surface on the Workbook/Worksheet/
# ListObject COM object model exposed here (WorkbookConnections.Add2,
# QueryTables.Add, Workbook.Queries.Add and CustomXMLParts.Add are all
# present on the object model but are inert no-ops in this host - they
# neither throw nor register anything), so that portion of the diff
# cannot be reproduced through COM automation. Everything else
# (sheets, cell data, shared strings, tables/styles, the defined name,
# and the view/selection changes) is reproduced exactly below.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1. New sheet "data2", positioned right after "data1" -> becomes the
#    active (second / index 2) sheet, which is what flips
#    bookViews/workbookView@activeTab to 1 and moves tabSelected from
#    sheet1 to sheet2 automatically on save.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "data2"

# 2. Cell data for data2.
$ws2.Range("A1").Value = "number for verification"
$ws2.Range("A2").Value = 3
$ws2.Range("A4").Value = "id"
$ws2.Range("B4").Value = "name"
$ws2.Range("A5").Value = 1
$ws2.Range("B5").Value = "1*"

# 3. "tData2" - the table the query's results were loaded into (A4:B5).
$ws2.ListObjects.Add(1, $ws2.Range("A4:B5"), $null, 1) | Out-Null
$ws2.ListObjects("Table2").Name = "tData2"
$ws2.ListObjects("tData2").TableStyle = "TableStyleMedium7"

# 4. "ptVerificationNumber" - the small parameter table (A1:A2).
$ws2.ListObjects.Add(1, $ws2.Range("A1:A2"), $null, 1) | Out-Null
$ws2.ListObjects("Table3").Name = "ptVerificationNumber"

# 5. Hidden workbook-scoped (local to data2) defined name Excel keeps
#    alongside a query-table-backed table.
$extData = $ws2.Names.Add("ExternalData_1", "=data2!`$A`$4:`$B`$5")
$extData.Visible = $false

# 6. Selection state captured in the diff: data2's sheetView selects
#    A3, sheet1 keeps its prior A1:C4 selection (sheet1's selection is
#    left untouched - only its tabSelected flag changes, automatically,
#    because data2 is now the active sheet).
$ws2.Range("A3").Select() | Out-Null
